# "Added more missing stars" - mark most checklist items as Done (B column
# checkbox = TRUE) and clear the old status/notes text that lived in column C,
# now that the sharedStrings table has been trimmed down to just the
# constellation names, "Done", and "Notes".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C cells whose old status note ("Missing", "Error", "Skewed", ...)
# is being removed entirely (no replacement text/checkbox).
$clearOnlyAddresses = "C31,C37,C39,C43,C61,C69,C73,C86"
foreach ($area in $ws.Range($clearOnlyAddresses).Areas) {
    $area.ClearContents()
}

# Row 63 loses its checked box (goes back to a blank row, like the
# clear-only rows above).
$ws.Range("B63").ClearContents()

# Column C cells whose old status note is removed AND replaced by ticking
# the "Done" checkbox in column B.
$clearAndCheckAddresses = "C2,C4,C5,C9,C11,C12,C17,C18,C24,C27,C40,C44,C46,C58,C59,C65,C67,C74,C75,C77,C79,C83,C84,C85,C87"
foreach ($area in $ws.Range($clearAndCheckAddresses).Areas) {
    $area.ClearContents()
}

$checkAddresses = "B2,B4,B5,B9,B11,B12,B17,B18,B24,B27,B40,B44,B46,B58,B59,B65,B67,B74,B75,B77,B79,B83,B84,B85,B87"
foreach ($area in $ws.Range($checkAddresses).Areas) {
    $area.Value = $true
}

# The header note in C1 still reads "Notes" - rewriting it lands on the new,
# compacted shared-string slot for the same text.
$ws.Range("C1").Value = "Notes"

# Restore the view: scrolled back to the top, with C31 selected (instead of
# the scrolled-down J51 selection saved before).
$ws.Range("A1").Select()
$ws.Range("C31").Select()
